# "fixed daily measurement file format issue"
#
# The "date" column (column E) on the "Metadata" sheet was stored as a
# garbled/truncated number (2031111 or 2023111) instead of the intended
# YYYYMMDD value 20231111 (2023-11-11). A block of rows for 2023-11-12
# (stored in "Box 3") were also missing their date (E) and Box (H)
# entries altogether, so those get filled in too.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Rows whose mis-typed date value needs correcting to 20231111.
$rowsToFixDate = @(143..220) + @(227..242) + @(247..250)
foreach ($r in $rowsToFixDate) {
    $ws.Cells.Item($r, 5).Value = 20231111
}

# Rows that were missing their date (E) and Box (H) entries entirely --
# these are the 2023-11-12 rows filed under "Box 3".
$rowsToAddCells = @(221..226) + @(243..246)
foreach ($r in $rowsToAddCells) {
    $ws.Cells.Item($r, 5).Value = 20231112
    $ws.Cells.Item($r, 8).Value = "Box 3"
}

# Leave the view roughly where the editor left it (scroll position +
# active cell) after making the fix.
$win = $excel.ActiveWindow
$win.ScrollRow = 216
$win.ScrollColumn = 1
$ws.Range("D233").Select()
